$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("incidentes")
$ws.Name = "Folha1"
$ws.Activate()
$ws.Range("G18").Select()
